# Qualification-Template.xlsx update (commit: "2 excel to quali (#19)")
#  * Rename sheet "GlobalAxes_DDI_PredVsObs" -> "GlobalAxesSettings"
#  * Add a new leading "Plot" column on that sheet (header row shifts right)
#  * Projects sheet: header "ID" -> "Id", make it the active/selected sheet
#  * Misc window/view bookkeeping to mirror the authoring session

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. GlobalAxes_DDI_PredVsObs -> GlobalAxesSettings, insert "Plot" column
# ---------------------------------------------------------------------
$axesSheet = $wb.Worksheets.Item("GlobalAxes_DDI_PredVsObs")

# Insert a new column A (existing columns shift right: A->B, B->C, ...)
$axesSheet.Columns.Item(1).Insert()

# New column header value
$axesSheet.Range("A1").Value = "Plot"

# Match the header formatting used by the rest of row 1 (bold + fill)
$axesSheet.Range("B1").Copy()
$axesSheet.Range("A1").PasteSpecial(-4122)
[void]($axesSheet.Application.CutCopyMode = $false)

# Column width for the freshly inserted column
$axesSheet.Columns.Item(1).ColumnWidth = 4.08984375

# View state: selection moves, sheet no longer the tab-selected one
[void]$axesSheet.Range("C6").Select()

# Rename the sheet itself
$axesSheet.Name = "GlobalAxesSettings"

# ---------------------------------------------------------------------
# 2. Projects sheet: "ID" header -> "Id", becomes the active sheet/tab
# ---------------------------------------------------------------------
$projectsSheet = $wb.Worksheets.Item("Projects")
$projectsSheet.Range("A1").Value = "Id"
[void]$projectsSheet.Activate()
[void]$projectsSheet.Range("A2").Select()

Write-Output "done"
